$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells we touch to Text format so numeric-looking
# strings (e.g. "0.999", "6.32") are not coerced into floating point numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated values scraped on Tue Jul 16 08:50:04 UTC 2024
$ws.Range("D2").Value = '62.721.63'
$ws.Range("E2").Value = '  -0.53%  '
$ws.Range("D3").Value = '3.364.07'
$ws.Range("E3").Value = '  -0.05%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '560.07'
$ws.Range("E5").Value = '  -0.49%  '
$ws.Range("D6").Value = '153.76'
$ws.Range("E6").Value = '  -0.13%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.23%  '
$ws.Range("D8").Value = '3.364.41'
$ws.Range("E8").Value = '  -0.19%  '
$ws.Range("E9").Value = '  +0.86%  '
$ws.Range("D10").Value = '7.36'
$ws.Range("E10").Value = '  -2.32%  '
$ws.Range("E11").Value = '  +0.40%  '
$ws.Range("E12").Value = '  -2.13%  '
$ws.Range("D13").Value = '3.945.37'
$ws.Range("E13").Value = '  +0.07%  '
$ws.Range("E14").Value = '  -3.86%  '
$ws.Range("E15").Value = '  +2.30%  '
$ws.Range("D16").Value = '26.83'
$ws.Range("E16").Value = '  -1.59%  '
$ws.Range("D17").Value = '62.794.99'
$ws.Range("E17").Value = '  -0.34%  '
$ws.Range("D18").Value = '3.288.13'
$ws.Range("E18").Value = '  -1.21%  '
$ws.Range("D19").Value = '6.19'
$ws.Range("E19").Value = '  -4.32%  '
$ws.Range("D20").Value = '13.94'
$ws.Range("E20").Value = '  +0.47%  '
$ws.Range("D21").Value = '372.88'
$ws.Range("E21").Value = '  -4.35%  '
$ws.Range("D22").Value = '7.96'
$ws.Range("E22").Value = '  -5.86%  '
$ws.Range("D23").Value = '0.995'
$ws.Range("E23").Value = '  -0.61%  '
$ws.Range("D24").Value = '70.81'
$ws.Range("E24").Value = '  +0.42%  '
$ws.Range("D25").Value = '0.524'
$ws.Range("E25").Value = '  -3.31%  '
$ws.Range("E26").Value = '  +16.30%  '
$ws.Range("D27").Value = '9.44'
$ws.Range("E27").Value = '  +6.41%  '
$ws.Range("E28").Value = '  -3.04%  '
$ws.Range("D29").Value = '1.01'
$ws.Range("E29").Value = '  +0.76%  '
$ws.Range("D30").Value = '6.03'
$ws.Range("E30").Value = '  +6.06%  '
$ws.Range("E31").Value = '  -2.47%  '
$ws.Range("E32").Value = '  +1.55%  '
$ws.Range("E33").Value = '  -1.85%  '
$ws.Range("D34").Value = '22.99'
$ws.Range("E34").Value = '  -0.43%  '
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("D36").Value = '6.70'
$ws.Range("E36").Value = '  -0.90%  '
$ws.Range("D37").Value = '159.11'
$ws.Range("E37").Value = '  -0.94%  '
$ws.Range("D38").Value = '1.44'
$ws.Range("E38").Value = '  -2.99%  '
$ws.Range("D39").Value = '0.0762'
$ws.Range("E39").Value = '  +1.86%  '
$ws.Range("D40").Value = '2.910.41'
$ws.Range("E40").Value = '  +2.31%  '
$ws.Range("D41").Value = '26.86'
$ws.Range("E41").Value = '  -0.68%  '
$ws.Range("E42").Value = '  -5.20%  '
$ws.Range("D43").Value = '0.0315'
$ws.Range("E43").Value = '  +0.29%  '
$ws.Range("D44").Value = '41.26'
$ws.Range("E44").Value = '  +1.03%  '
$ws.Range("D45").Value = '4.28'
$ws.Range("E45").Value = '  -1.01%  '
$ws.Range("D46").Value = '0.739'
$ws.Range("E46").Value = '  -1.54%  '
$ws.Range("D47").Value = '22.88'
$ws.Range("E47").Value = '  +2.78%  '
$ws.Range("E48").Value = '  +0.16%  '
$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").Value = '6.32'
$ws.Range("E49").Value = '  -0.02%  '
$ws.Range("B50").Value = 'dogwifhat'
$ws.Range("C50").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D50").Value = '2.08'
$ws.Range("E50").Value = '  +14.01%  '
$ws.Range("D51").Value = '0.824'
$ws.Range("E51").Value = '  +1.61%  '
